$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 124
$ws.Range("I12").Value = 124
$ws.Range("K12").Value = 124
$ws.Range("M12").Value = 46

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 10000
$ws.Range("I21").Value = 10000
$ws.Range("K21").Value = 10000
$ws.Range("M21").Value = -9532

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 10000
$ws.Range("I23").Value = 10000
$ws.Range("K23").Value = 10000
$ws.Range("M23").Value = -9766

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 984.5
$ws.Range("I41").Value = 972.5
$ws.Range("J41").Value = 996.5
$ws.Range("K41").Value = 972.5
$ws.Range("L41").Value = 996.5
$ws.Range("M41").Value = -532.5
$ws.Range("N41").Value = -1876.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 2962.6667
$ws.Range("J48").Value = 2962.6667
$ws.Range("L48").Value = 8888.000100000001
$ws.Range("N48").Value = -9472.000100000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H56").Value = 2962.6667
$ws.Range("J56").Value = 2962.6667
$ws.Range("L56").Value = 8888.000100000001
$ws.Range("N56").Value = -9956.000100000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 884.5
$ws.Range("I99").Value = 289.5
$ws.Range("J99").Value = 1479.5
$ws.Range("K99").Value = 868.5
$ws.Range("L99").Value = 4438.5
$ws.Range("M99").Value = 629.5
$ws.Range("N99").Value = -7434.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 33333650
$ws.Range("I101").Value = 50000124
$ws.Range("J101").Value = 700
$ws.Range("K101").Value = 150000372
$ws.Range("L101").Value = 2100
$ws.Range("M101").Value = -149998750
$ws.Range("N101").Value = -5344

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1299.6
$ws.Range("I111").Value = 1350
$ws.Range("K111").Value = 4050
$ws.Range("M111").Value = -983

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3985.5715
$ws.Range("I132").Value = 2570.647
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 7711.941
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -5181.941
$ws.Range("N132").Value = -35057

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 3
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2444.6765
$ws.Range("I32").Value = 2510.4849
$ws.Range("K32").Value = 2510.4849
$ws.Range("M32").Value = -2223.4849

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 3999.5
$ws.Range("J46").Value = 3999
$ws.Range("L46").Value = 3999
$ws.Range("N46").Value = -4637

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2204.0908
$ws.Range("I74").Value = 1809.375
$ws.Range("J74").Value = 3256.6667
$ws.Range("K74").Value = 1809.375
$ws.Range("L74").Value = 3256.6667
$ws.Range("M74").Value = -935.375
$ws.Range("N74").Value = -5004.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2204.0908
$ws.Range("I77").Value = 1809.375
$ws.Range("J77").Value = 3256.6667
$ws.Range("K77").Value = 9046.875
$ws.Range("L77").Value = 16283.3335
$ws.Range("M77").Value = -4678.875
$ws.Range("N77").Value = -25019.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 12333999
$ws.Range("I110").Value = 12333999
$ws.Range("K110").Value = 12333999
$ws.Range("M110").Value = -12331954

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3886.375
$ws.Range("I122").Value = 2700
$ws.Range("J122").Value = 4281.8335
$ws.Range("K122").Value = 8100
$ws.Range("L122").Value = 12845.5005
$ws.Range("M122").Value = -5650
$ws.Range("N122").Value = -17745.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5430.0835
$ws.Range("I107").Value = 4183
$ws.Range("K107").Value = 4183
$ws.Range("M107").Value = -2263

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3109.9443
$ws.Range("I31").Value = 3092.625
$ws.Range("J31").Value = 3248.5
$ws.Range("K31").Value = 3092.625
$ws.Range("L31").Value = 3248.5
$ws.Range("M31").Value = -2797.625
$ws.Range("N31").Value = -3838.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3109.9443
$ws.Range("I34").Value = 3092.625
$ws.Range("J34").Value = 3248.5
$ws.Range("K34").Value = 3092.625
$ws.Range("L34").Value = 3248.5
$ws.Range("M34").Value = -2890.625
$ws.Range("N34").Value = -3652.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5604.4287
$ws.Range("I99").Value = 1033
$ws.Range("K99").Value = 1033
$ws.Range("M99").Value = 465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3649.4614
$ws.Range("I105").Value = 3157.3333
$ws.Range("K105").Value = 3157.3333
$ws.Range("M105").Value = -1410.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1228.4667
$ws.Range("I107").Value = 1054.25
$ws.Range("K107").Value = 1054.25
$ws.Range("M107").Value = 865.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5604.4287
$ws.Range("I126").Value = 1033
$ws.Range("K126").Value = 3099
$ws.Range("M126").Value = -629

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2700.95
$ws.Range("I132").Value = 2940.2727
$ws.Range("J132").Value = 2408.4443
$ws.Range("K132").Value = 8820.8181
$ws.Range("L132").Value = 7225.3329
$ws.Range("M132").Value = -6290.8181
$ws.Range("N132").Value = -12285.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 124.9
$ws.Range("I26").Value = 93.625
$ws.Range("K26").Value = 280.875
$ws.Range("M26").Value = 7.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 14506.667
$ws.Range("I75").Value = 2013
$ws.Range("J75").Value = 17005.4
$ws.Range("K75").Value = 6039
$ws.Range("L75").Value = 51016.2
$ws.Range("M75").Value = -5041
$ws.Range("N75").Value = -53012.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 14506.667
$ws.Range("I78").Value = 2013
$ws.Range("J78").Value = 17005.4
$ws.Range("K78").Value = 18117
$ws.Range("L78").Value = 153048.6
$ws.Range("M78").Value = -13125
$ws.Range("N78").Value = -163032.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1671.8182
$ws.Range("J113").Value = 1599.4445
$ws.Range("L113").Value = 4798.333500000001
$ws.Range("N113").Value = -9138.333500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3175.182
$ws.Range("I132").Value = 1673.1666
$ws.Range("K132").Value = 15058.4994
$ws.Range("M132").Value = -12528.4994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5971.5557
$ws.Range("I126").Value = 6427.7144
$ws.Range("J126").Value = 4375
$ws.Range("K126").Value = 19283.1432
$ws.Range("L126").Value = 13125
$ws.Range("M126").Value = -16813.1432
$ws.Range("N126").Value = -18065

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 37499.5
$ws.Range("J134").Value = 37499.5
$ws.Range("L134").Value = 112498.5
$ws.Range("N134").Value = -117568.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5463.28
$ws.Range("I7").Value = 2008.5454
$ws.Range("J7").Value = 8177.7144
$ws.Range("K7").Value = 2008.5454
$ws.Range("L7").Value = 8177.7144
$ws.Range("M7").Value = -1896.5454
$ws.Range("N7").Value = -8401.714400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2916.8823
$ws.Range("I40").Value = 2353
$ws.Range("J40").Value = 4749.5
$ws.Range("K40").Value = 2353
$ws.Range("L40").Value = 4749.5
$ws.Range("M40").Value = -2217
$ws.Range("N40").Value = -5021.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1388.7
$ws.Range("J82").Value = 963.3333
$ws.Range("L82").Value = 963.3333
$ws.Range("N82").Value = -1685.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1388.7
$ws.Range("J85").Value = 963.3333
$ws.Range("L85").Value = 963.3333
$ws.Range("N85").Value = -3459.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5854.0645
$ws.Range("I122").Value = 5656.885
$ws.Range("K122").Value = 16970.655
$ws.Range("M122").Value = -14520.655

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5463.28
$ws.Range("I126").Value = 2008.5454
$ws.Range("J126").Value = 8177.7144
$ws.Range("K126").Value = 6025.6362
$ws.Range("L126").Value = 24533.1432
$ws.Range("M126").Value = -3555.6362
$ws.Range("N126").Value = -29473.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2239.4
$ws.Range("I132").Value = 1699.25
$ws.Range("K132").Value = 5097.75
$ws.Range("M132").Value = -2567.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 6251037
$ws.Range("I100").Value = 9092015
$ws.Range("K100").Value = 18184030
$ws.Range("M100").Value = -18183489

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1853.8889
$ws.Range("I122").Value = 1470.8572
$ws.Range("J122").Value = 3194.5
$ws.Range("K122").Value = 4412.571599999999
$ws.Range("L122").Value = 9583.5
$ws.Range("M122").Value = -1962.571599999999
$ws.Range("N122").Value = -14483.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8119.75
$ws.Range("I132").Value = 7994.4287
$ws.Range("J132").Value = 8997
$ws.Range("K132").Value = 23983.2861
$ws.Range("L132").Value = 26991
$ws.Range("M132").Value = -21453.2861
$ws.Range("N132").Value = -32051
